# Updated symbol list (crypto price refresh) — column D ("Price") values
# for the affected rows. Source cells are plain text (e.g. "243.64"),
# so force a text number-format before assigning the new string value —
# this keeps the cell stored as text (preserving things like trailing
# zeros: "5.300", "0.001539", etc.) instead of Excel auto-coercing the
# numeric-looking string into a numeric cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = $origStyle
}

Set-TextValue "D2"  "243.53"
Set-TextValue "D3"  "23.52"
Set-TextValue "D4"  "5.296"
Set-TextValue "D5"  "0.05768"
Set-TextValue "D6"  "6.470"
Set-TextValue "D7"  "3.336"
Set-TextValue "D8"  "0.8109"
Set-TextValue "D9"  "0.8799"
Set-TextValue "D10" "0.1379"
Set-TextValue "D11" "0.07296"
Set-TextValue "D12" "0.03086"
Set-TextValue "D13" "0.03058"
Set-TextValue "D15" "3.858"
Set-TextValue "D16" "0.001547"
Set-TextValue "D17" "0.04724"
Set-TextValue "D18" "0.0006069"
Set-TextValue "D19" "0.006036"
Set-TextValue "D20" "0.001296"
Set-TextValue "D21" "0.004602"
Set-TextValue "D22" "0.00008807"
Set-TextValue "D24" "2.141"
Set-TextValue "D25" "0.3182"
Set-TextValue "D41" "0.006445"
Set-TextValue "D42" "0.1053"
Set-TextValue "D44" "0.008023"
Set-TextValue "D45" "0.00005479"
Set-TextValue "D47" "0.6009"
Set-TextValue "D48" "0.001860"
